$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value (serial 45189) for every data
# row (2 through 517). Update it to the new serial date value 45190.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 517
}

$ws.Range("C2:C$lastRow").Value = 45190
